# Update "detection field test data.xlsx" - append new field test rows
# for Britt/Koda and Sean/Koda sessions on 2025-06-07 to the "human" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("human")

# The Date (col B) and Time of search (col C) columns carry existing
# number-formatted styles (short date / time-of-day). Copy the formatting
# down from the row above before writing values so the new rows reuse the
# workbook's existing style entries instead of creating new ones.
$ws.Cells.Item(11, 2).Copy($ws.Cells.Item(12, 2))
$ws.Cells.Item(11, 3).Copy($ws.Cells.Item(12, 3))
$ws.Cells.Item(11, 2).Copy($ws.Cells.Item(13, 2))
$ws.Cells.Item(11, 3).Copy($ws.Cells.Item(13, 3))
$ws.Cells.Item(11, 2).Copy($ws.Cells.Item(14, 2))
$ws.Cells.Item(11, 3).Copy($ws.Cells.Item(14, 3))
$ws.Cells.Item(11, 2).Copy($ws.Cells.Item(15, 2))
$ws.Cells.Item(11, 3).Copy($ws.Cells.Item(15, 3))

# Row 12 - Britt, target found
$ws.Cells.Item(12, 1).Value = "Britt"
$ws.Cells.Item(12, 2).Value = 45815
$ws.Cells.Item(12, 3).Value = 0.42708333333333331
$ws.Cells.Item(12, 4).Value = "Overcast, cool."
$ws.Cells.Item(12, 5).Value = $true
$ws.Cells.Item(12, 6).Value = "32 minutes 30 seconds"
$ws.Cells.Item(12, 7).Formula = "=32*60+30"
$ws.Cells.Item(12, 8).Value = "Transects (short length). A few minutes of rain."

# Row 13 - Britt, not found (1 hour elapsed)
$ws.Cells.Item(13, 1).Value = "Britt"
$ws.Cells.Item(13, 2).Value = 45815
$ws.Cells.Item(13, 3).Value = 0.47916666666666669
$ws.Cells.Item(13, 4).Value = "Overcast, cool."
$ws.Cells.Item(13, 5).Value = $false
$ws.Cells.Item(13, 6).Value = "1 hour"
$ws.Cells.Item(13, 7).Value = 3600
$ws.Cells.Item(13, 8).Value = "Transects (short length). Not found and funnily I couldn't refind either, had to get Koda to retrieve the placement."

# Row 14 - Sean, target found
$ws.Cells.Item(14, 1).Value = "Sean"
$ws.Cells.Item(14, 2).Value = 45815
$ws.Cells.Item(14, 3).Value = 0.58333333333333337
$ws.Cells.Item(14, 4).Value = "Overcast, cold."
$ws.Cells.Item(14, 5).Value = $true
$ws.Cells.Item(14, 6).Value = "27 minutes 1 second"
$ws.Cells.Item(14, 7).Value = 1621
$ws.Cells.Item(14, 8).Value = "Random walk then transects. Found on first intensive transect."

# Row 15 - Sean, not found (1 hour elapsed)
$ws.Cells.Item(15, 1).Value = "Sean"
$ws.Cells.Item(15, 2).Value = 45815
$ws.Cells.Item(15, 3).Value = 0.61458333333333337
$ws.Cells.Item(15, 4).Value = "Overcast, cold."
$ws.Cells.Item(15, 5).Value = $false
$ws.Cells.Item(15, 6).Value = "1 hour"
$ws.Cells.Item(15, 7).Value = 3600
$ws.Cells.Item(15, 8).Value = "Random walk then transects. Was very cold."

$ws.Range("G24").Select()
